# Applies two edits to JourneyPPT1.MrigaArora.pptx:
#  1. Slide 14 ("Content Placeholder 2"): merge the two runs of the
#     "Solution: ..." paragraph into a single run with the combined text.
#  2. Slide 6 ("Content Placeholder 2"): bold + single-underline the whole
#     "Technical Boot Camp Day 4 and Day 5: ..." paragraph (all 5 runs).

$p = $ppt.ActivePresentation

# --- Edit 1: slide 14, merge "Solution: ... " + "understand." runs ---
$slide14 = $p.Slides.Item(14)
$shape14 = $slide14.Shapes.Item(5)   # "Content Placeholder 2"
$tr14 = $shape14.TextFrame.TextRange

$para2 = $tr14.Paragraphs(2, 1)

# Paragraphs().Text/.Length include the trailing paragraph-mark (CR)
# character for every paragraph except the very last one in the text
# frame, so strip it before doing character-offset math.
$para2Text = $para2.Text
if ($para2Text.Length -gt 0 -and [int]$para2Text[$para2Text.Length - 1] -eq 13) {
    $para2Text = $para2Text.Substring(0, $para2Text.Length - 1)
}

$run1 = $para2.Runs(1, 1)
$run1Len = $run1.Text.Length
$remainderStart = $para2.Start + $run1Len
$remainderLen = $para2Text.Length - $run1Len
$run2 = $tr14.Characters($remainderStart, $remainderLen)

$tailText = $run2.Text
$run2.Delete()
$run1.Text = $run1.Text + $tailText

# --- Edit 2: slide 6, bold + underline the "Technical Boot Camp" heading ---
$slide6 = $p.Slides.Item(6)
$shape6 = $slide6.Shapes.Item(6)   # "Content Placeholder 2"
$tr6 = $shape6.TextFrame.TextRange

$para1 = $tr6.Paragraphs(1, 1)
$runs1 = $para1.Runs()
for ($i = 1; $i -le $runs1.Count; $i++) {
    $run = $para1.Runs($i, 1)
    $run.Font.Bold = -1
    $run.Font.Underline = -1
}
